$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Clear the cells whose string values are being replaced, so the shared
# strings table drops the old, now-unused entries before new ones are
# appended (mirrors how Excel compacts/rebuilds the table on save).
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("B1").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()

# Re-enter the new values in the order that reproduces the target shared
# string table ordering: TEST1, TEST2, TEST3, Name Detail, TEST
$ws.Range("B2").Value = "TEST1"
$ws.Range("B3").Value = "TEST2"
$ws.Range("B4").Value = "TEST3"
$ws.Range("B1").Value = "Name Detail"
$ws.Range("A2").Value = "TEST"
$ws.Range("A3").Value = "TEST"
$ws.Range("A4").Value = "TEST"

$ws.Range("A2:A4").Select()
$ws.Application.ActiveCell = $ws.Range("A2")
